# Add a "Db ID" column to the Oracle DB export template.
#
# The new column is inserted immediately before the "Hostname" column
# (which is column C, since A = "DB Name" and B = "Unique name"), pushing
# "Hostname" through "Clusterware Version" one column to the right.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at C (before "Hostname"); Excel shifts
# C:V -> D:W and copies column C's formatting from its left neighbor (B).
$ws.Columns.Item(3).Insert()
$ws.Cells.Item(1, 3).Value = "Db ID"

# The header row uses one consistent look: bold Calibri 11, centered,
# no wrap. Re-apply it across the whole (now 23-column) header row so the
# freshly inserted "Db ID" header matches its neighbors exactly.
$headerRow = $ws.Range("A1:W1")
$headerRow.Font.Name = "Calibri"
$headerRow.Font.Size = 11
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108
$headerRow.VerticalAlignment = -4108
$headerRow.WrapText = $false

# "Pgsql Migrability" (shifted from U1 to V1 by the insert) is the one
# header that wraps its text onto two lines.
$ws.Range("V1").WrapText = $true

# Leave the selection where the editor last left it.
$ws.Range("C9").Select()
